$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 test data: browser IE first, then the new (longer) test case
# name, then re-assert the Automation user. This ordering keeps the shared
# string table aligned with how the workbook was actually re-saved.
$ws.Range("B2").Value = "IE"
$ws.Range("A2").Value = "Exchange_3_Return_To_Appointments_List_Page_From_Appointment_Details_Page_IE"
$ws.Range("C2").Value = "Automation"

# Widen column A to fit the new (longer) test case name
$ws.Columns.Item(1).ColumnWidth = 72.83072916666667

# Move the active selection to A6
[void]$ws.Range("A6").Select()
